$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1399.1052
$ws.Range("I43").Value = 1055.1666
$ws.Range("J43").Value = 1557.8462
$ws.Range("K43").Value = 1055.1666
$ws.Range("L43").Value = 1557.8462
$ws.Range("M43").Value = -986.1666
$ws.Range("N43").Value = -1695.8462

$ws.Range("H98").Value = 2270.476
$ws.Range("I98").Value = 2090.9
$ws.Range("J98").Value = 2433.7273
$ws.Range("K98").Value = 2090.9
$ws.Range("L98").Value = 2433.7273
$ws.Range("M98").Value = -592.9000000000001
$ws.Range("N98").Value = -5429.7273

$ws.Range("H106").Value = 3364.75
$ws.Range("I106").Value = 3297.125
$ws.Range("J106").Value = 3500
$ws.Range("K106").Value = 3297.125
$ws.Range("L106").Value = 3500
$ws.Range("M106").Value = -2666.125
$ws.Range("N106").Value = -4762

$ws.Range("H122").Value = 2270.476
$ws.Range("I122").Value = 2090.9
$ws.Range("J122").Value = 2433.7273
$ws.Range("K122").Value = 6272.700000000001
$ws.Range("L122").Value = 7301.1819
$ws.Range("M122").Value = -3822.700000000001
$ws.Range("N122").Value = -12201.1819

$ws.Range("H123").Value = 25373.334
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 25373.334
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 25373.334
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -35173.334

$ws.Range("H126").Value = 40780
$ws.Range("J126").Value = 40780
$ws.Range("L126").Value = 40780
$ws.Range("N126").Value = -50660

$ws.Range("H130").Value = 97528.57000000001
$ws.Range("J130").Value = 97528.57000000001
$ws.Range("L130").Value = 97528.57000000001
$ws.Range("N130").Value = -107568.57

$ws.Range("H134").Value = 42727.273
$ws.Range("J134").Value = 42727.273
$ws.Range("L134").Value = 42727.273
$ws.Range("N134").Value = -52867.273

$ws.Range("H140").Value = 38005.266
$ws.Range("J140").Value = 38005.266
$ws.Range("L140").Value = 38005.266
$ws.Range("N140").Value = -48365.266

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6024.77
$ws.Range("I32").Value = 4457.554
$ws.Range("J32").Value = 13676.471
$ws.Range("K32").Value = 4457.554
$ws.Range("L32").Value = 13676.471
$ws.Range("M32").Value = -4170.554
$ws.Range("N32").Value = -14250.471

$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 30000
$ws.Range("N134").Value = -40140

$ws.Range("H135").Value = 19500.889
$ws.Range("J135").Value = 19500.889
$ws.Range("L135").Value = 19500.889
$ws.Range("N135").Value = -29640.889

$ws.Range("H139").Value = 31714.428
$ws.Range("J139").Value = 31714.428
$ws.Range("L139").Value = 31714.428
$ws.Range("N139").Value = -41994.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 35780.47
$ws.Range("J135").Value = 37534.938
$ws.Range("L135").Value = 37534.938
$ws.Range("N135").Value = -47674.938

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 17384.6
$ws.Range("J28").Value = 17384.6
$ws.Range("L28").Value = 17384.6
$ws.Range("N28").Value = -17874.6

$ws.Range("H100").Value = 40779.5
$ws.Range("J100").Value = 40779.5
$ws.Range("L100").Value = 40779.5
$ws.Range("N100").Value = -42943.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3628.4614
$ws.Range("I134").Value = 1154.3478
$ws.Range("J134").Value = 4983.3335
$ws.Range("K134").Value = 3463.0434
$ws.Range("L134").Value = 14950.0005
$ws.Range("M134").Value = 1606.9566
$ws.Range("N134").Value = -25090.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1964.24
$ws.Range("I122").Value = 1732.1177
$ws.Range("J122").Value = 2457.5
$ws.Range("K122").Value = 5196.3531
$ws.Range("L122").Value = 7372.5
$ws.Range("M122").Value = -2746.3531
$ws.Range("N122").Value = -12272.5

$ws.Range("H141").Value = 46574.75
$ws.Range("J141").Value = 46574.75
$ws.Range("L141").Value = 46574.75
$ws.Range("N141").Value = -56934.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 166669020
$ws.Range("I7").Value = 500000740
$ws.Range("J7").Value = 3151.25
$ws.Range("K7").Value = 500000740
$ws.Range("L7").Value = 3151.25
$ws.Range("M7").Value = -500000628
$ws.Range("N7").Value = -3375.25

$ws.Range("H40").Value = 2517.4614
$ws.Range("I40").Value = 1967.6923
$ws.Range("K40").Value = 1967.6923
$ws.Range("M40").Value = -1831.6923

$ws.Range("H122").Value = 3371.0386
$ws.Range("I122").Value = 3292.7144
$ws.Range("J122").Value = 3700
$ws.Range("K122").Value = 9878.143199999999
$ws.Range("L122").Value = 11100
$ws.Range("M122").Value = -7428.143199999999
$ws.Range("N122").Value = -16000

$ws.Range("H126").Value = 166669020
$ws.Range("I126").Value = 500000740
$ws.Range("J126").Value = 3151.25
$ws.Range("K126").Value = 1500002220
$ws.Range("L126").Value = 9453.75
$ws.Range("M126").Value = -1499999750
$ws.Range("N126").Value = -14393.75

$ws.Range("H127").Value = 27700
$ws.Range("J127").Value = 27700
$ws.Range("L127").Value = 27700
$ws.Range("N127").Value = -37620

$ws.Range("H130").Value = 24894.445
$ws.Range("J130").Value = 21610
$ws.Range("L130").Value = 21610
$ws.Range("N130").Value = -31650

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 30638.125
$ws.Range("I70").Value = 20000
$ws.Range("J70").Value = 32157.857
$ws.Range("K70").Value = 20000
$ws.Range("L70").Value = 32157.857
$ws.Range("M70").Value = -19685
$ws.Range("N70").Value = -32787.857

$ws.Range("H73").Value = 30638.125
$ws.Range("I73").Value = 20000
$ws.Range("J73").Value = 32157.857
$ws.Range("K73").Value = 20000
$ws.Range("L73").Value = 32157.857
$ws.Range("M73").Value = -18908
$ws.Range("N73").Value = -34341.857

$ws.Range("H107").Value = 449.9091
$ws.Range("I107").Value = 362.66666
$ws.Range("J107").Value = 554.6
$ws.Range("K107").Value = 1087.99998
$ws.Range("L107").Value = 1663.8
$ws.Range("M107").Value = 832.0000199999999
$ws.Range("N107").Value = -5503.8

$ws.Range("H122").Value = 1230.3103
$ws.Range("I122").Value = 854.95
$ws.Range("J122").Value = 2064.4443
$ws.Range("K122").Value = 2564.85
$ws.Range("L122").Value = 6193.3329
$ws.Range("M122").Value = -114.8500000000004
$ws.Range("N122").Value = -11093.3329

$ws.Range("H123").Value = 29756
$ws.Range("J123").Value = 29756
$ws.Range("L123").Value = 29756
$ws.Range("N123").Value = -39556

$ws.Range("H125").Value = 43620.312
$ws.Range("J125").Value = 43620.312
$ws.Range("L125").Value = 43620.312
$ws.Range("N125").Value = -53460.312

$ws.Range("H135").Value = 32400
$ws.Range("J135").Value = 32400
$ws.Range("L135").Value = 32400
$ws.Range("N135").Value = -42540

$ws.Range("H140").Value = 69561
$ws.Range("J140").Value = 69561
$ws.Range("L140").Value = 69561
$ws.Range("N140").Value = -79921

$ws.Range("H141").Value = 57893.168
$ws.Range("J141").Value = 57893.168
$ws.Range("L141").Value = 57893.168
$ws.Range("N141").Value = -68253.16800000001
